# updated topics in syllabus
# Rewrite five topic cells in column H of the schedule sheet with the new
# curriculum topics. Cells are set in this order so that the shared-string
# table's new entries land at the same indices the saved workbook expects.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H47").Value = "Sorting, Exception handling"
$ws.Range("H43").Value = "Random Module, Data Structures (tuples)"
$ws.Range("H33").Value = "Operators and Expressions"
$ws.Range("H34").Value = "Functions"
$ws.Range("H39").Value = "For loops, 2D lists"

# Move the active selection to where the author last left off editing.
$ws.Range("F47").Select()
